$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived data: rows with Target cluster = "ECs" have been dropped,
# and all downstream values recomputed. Columns: A Sending cluster,
# B Ligand symbol, C Receptor symbol, D Target cluster, E..T metrics.

$data = @(
    @("ECs",   "FAPs",  3, 1, 6.156604333333333, 18.469813,  0.3861700262161295, 0.3861700262161295, 3, 1,                  0.220618,            0.6618539999999999, 0.6621850925462731, 0.6621850925462731, 1.358257734811333,  12.224319613302,    0.2557160345485244,  0.2557160345485244),
    @("ECs",   "MuSCs", 3, 1, 6.156604333333333, 18.469813,  0.3861700262161295, 0.3861700262161295, 2, 0.6666666666666666, 0.1125486666666667, 0.337646,            0.3378149074537269, 0.3378149074537269, 0.6929176089108889, 6.236258480198,     0.1304539916676051,  0.1304539916676051),
    @("FAPs",  "FAPs",  3, 1, 5.867977666666667, 17.603933,  0.3680660582820729, 0.3680660582820729, 3, 1,                  0.220618,            0.6618539999999999, 0.6621850925462731, 0.6621850925462731, 1.294581496864667,  11.651233471782,    0.2437278568666564,  0.2437278568666564),
    @("FAPs",  "MuSCs", 3, 1, 5.867977666666667, 17.603933,  0.3680660582820729, 0.3680660582820729, 2, 0.6666666666666666, 0.1125486666666667, 0.337646,            0.3378149074537269, 0.3378149074537269, 0.6604330624131112, 5.943897561718001,  0.1243382014154165,  0.1243382014154165),
    @("MuSCs", "FAPs",  3, 1, 3.918147666666667, 11.754443,  0.2457639155017975, 0.2457639155017975, 3, 1,                  0.220618,            0.6618539999999999, 0.6621850925462731, 0.6621850925462731, 0.8644139019246667, 7.779725117322001,  0.1627412011310922,  0.1627412011310922),
    @("MuSCs", "MuSCs", 3, 1, 3.918147666666667, 11.754443,  0.2457639155017975, 0.2457639155017975, 2, 0.6666666666666666, 0.1125486666666667, 0.337646,            0.3378149074537269, 0.3378149074537269, 0.4409822956864445, 3.968840661178,     0.08302271437070528, 0.08302271437070528)
)

$row = 2
foreach ($r in $data) {
    $ws.Range("A$row").Value = $r[0]
    $ws.Range("B$row").Value = "Ntf3"
    $ws.Range("C$row").Value = "Ntrk1"
    $ws.Range("D$row").Value = $r[1]
    $ws.Range("E$row").Value = $r[2]
    $ws.Range("F$row").Value = $r[3]
    $ws.Range("G$row").Value = $r[4]
    $ws.Range("H$row").Value = $r[5]
    $ws.Range("I$row").Value = $r[6]
    $ws.Range("J$row").Value = $r[7]
    $ws.Range("K$row").Value = $r[8]
    $ws.Range("L$row").Value = $r[9]
    $ws.Range("M$row").Value = $r[10]
    $ws.Range("N$row").Value = $r[11]
    $ws.Range("O$row").Value = $r[12]
    $ws.Range("P$row").Value = $r[13]
    $ws.Range("Q$row").Value = $r[14]
    $ws.Range("R$row").Value = $r[15]
    $ws.Range("S$row").Value = $r[16]
    $ws.Range("T$row").Value = $r[17]
    $row++
}

# Rows 8:10 from the old (9-row) dataset are no longer present.
$ws.Range("A8:T10").Delete()
